# feat(filters): Add `merge` option in `dict` filter.
#
# Updates the "origin" sheet with a new recursive-dict example (A1 formula
# text gains a "merge":true option, plus new rows 3-4 demonstrating the
# `merge-1` / `merge-2` values), and references that new block from a new
# row on the "ref" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("origin")
$ws2 = $wb.Worksheets.Item("ref")

# "origin" sheet (sheet1.xml): update A1's filter text to include the new
# "merge" option, and add the merge-1 / merge-2 demo rows.
$ws1.Range("A1").Value = '#ref!A1(RD):._:R[{"fun":"dict","key":"upper","value":"ref","merge":true}]'

$ws1.Range("A3").Value = "merge-1"
$ws1.Range("B3").Value = 4
$ws1.Range("A4").Value = "merge-2"
$ws1.Range("B4").Value = "#B1"

# "ref" sheet (sheet2.xml): add a new row referencing the new origin block.
$ws2.Range("B25").Value = '#origin!A3:B4["recursive", "dict"]'

# Match the author's final selection/active-sheet state: "origin" becomes
# the active tab with C34 selected, "ref" keeps B25 selected.
[void]$ws2.Range("B25").Select()
[void]$ws1.Activate()
[void]$ws1.Range("C34").Select()
